# Automatische test-sync: 2025-06-19 21:51:50
# Append the newly received mail to the "Logs" sheet and refresh the
# "Dashboard" category-count summary to reflect it.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append new row 31 -----------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 31
$logs.Cells.Item($newRow, 1).Value = "Offerte voor zakelijke samenwerking"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Cells.Item($newRow, 4).Value = "Offerte / Prijsaanvraag"
# Column E (Antwoord) intentionally left blank - no reply sent yet.
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 21:51:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# ---- Extend the conditional formatting ranges to include the new row --
$catFormats = $logs.Range("D2:D30").FormatConditions
$catFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D31"))

$repliedFormats = $logs.Range("G2:G30").FormatConditions
$repliedFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G31"))

# ---- Dashboard sheet: refresh the category totals ---------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

# "Offerte / Prijsaanvraag" now has 3 occurrences (was 2) and moves above
# "Klacht / Probleem" (still 2) in the descending-count ordering.
$dashboard.Cells.Item(7, 1).Value = "Offerte / Prijsaanvraag"
$dashboard.Cells.Item(7, 2).Value = 3
$dashboard.Cells.Item(8, 1).Value = "Klacht / Probleem"
$dashboard.Cells.Item(8, 2).Value = 2
